$d = $word.ActiveDocument

# 1) Merge split runs for method-spec headings into single runs (removes gramStart/gramEnd proofErr marks)
$d.Content.Find.Execute("Đặc tả phương thức SuaThongTinTaiKhoan()", $true, $false, $false, $false, $false, $true, 1, $false, "Đặc tả phương thức SuaThongTinTaiKhoan()", 2) | Out-Null
$d.Content.Find.Execute("Đặc tả phương thức DangNhap()", $true, $false, $false, $false, $false, $true, 1, $false, "Đặc tả phương thức DangNhap()", 2) | Out-Null
$d.Content.Find.Execute("Đặc tả phương thức QuenMatKhau()", $true, $false, $false, $false, $false, $true, 1, $false, "Đặc tả phương thức QuenMatKhau()", 2) | Out-Null

# 2) Merge the "Nếu email không hợp lệ..." sentence into a single run as well
$emailMsg = "Nếu email không hợp lệ thì kết thúc và  hiển thị thông báo “Email không hợp lệ, kiểm tra lại”"
$d.Content.Find.Execute($emailMsg, $true, $false, $false, $false, $false, $true, 1, $false, $emailMsg, 2) | Out-Null

# 3) Append 4 blank paragraphs (class-diagram placeholder space) right after the last "Notes" section title
#    and before the final trailing blank paragraph of the document body.
$lastTable = $d.Tables($d.Tables.Count)
$afterTableRange = $d.Range($lastTable.Range.End, $d.Content.End)
$finalPara = $afterTableRange.Paragraphs(2)
$insertionPoint = $finalPara.Range
$insertionPoint.Collapse(1)
$blankParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr></w:p>' 
$fourBlankParas = $blankParaXml + $blankParaXml + $blankParaXml + $blankParaXml
$null = $insertionPoint.InsertXML($fourBlankParas)

Write-Host "edit applied"
